# Ajout de nouvelles espèces
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-strings table is rebuilt in first-use order, so we must set the
# new species names in the exact order their <si> entries appear in the
# target sharedStrings.xml (indices 42-50), not simply in row order.
# Row 43 reuses the existing "Dragon de Komodo" string (index 9), so it is
# set last among these new cells to avoid creating a spurious new entry.

$ws.Cells.Item(41, 1).Value = "Python royal"
$ws.Cells.Item(45, 1).Value = "Tortue léopard"
$ws.Cells.Item(49, 1).Value = "Scorpion"
$ws.Cells.Item(50, 1).Value = "Mante religieuse"
$ws.Cells.Item(42, 1).Value = "Dragon barbu de l'Est"
$ws.Cells.Item(44, 1).Value = "Caméléon panthère"
$ws.Cells.Item(47, 1).Value = "Sonneur oritental"
$ws.Cells.Item(48, 1).Value = "Tarentule vraie"
$ws.Cells.Item(46, 1).Value = "Rainette arboricole"
$ws.Cells.Item(43, 1).Value = "Dragon de Komodo"

# Numeric columns (diets_id / status_id) for the new rows.
$ws.Cells.Item(41, 2).Value = 8
$ws.Cells.Item(41, 3).Value = 7

$ws.Cells.Item(42, 2).Value = 6
$ws.Cells.Item(42, 3).Value = 7

$ws.Cells.Item(43, 2).Value = 8
$ws.Cells.Item(43, 3).Value = 5

$ws.Cells.Item(44, 2).Value = 6
$ws.Cells.Item(44, 3).Value = 7

$ws.Cells.Item(45, 2).Value = 8
$ws.Cells.Item(45, 3).Value = 7

$ws.Cells.Item(46, 2).Value = 6
$ws.Cells.Item(46, 3).Value = 7

$ws.Cells.Item(47, 2).Value = 6
$ws.Cells.Item(47, 3).Value = 7

$ws.Cells.Item(48, 2).Value = 6
$ws.Cells.Item(48, 3).Value = 7

$ws.Cells.Item(49, 2).Value = 8
$ws.Cells.Item(49, 3).Value = 7

$ws.Cells.Item(50, 2).Value = 6
$ws.Cells.Item(50, 3).Value = 7

# Update the view's current selection to reflect where the author ended up
# after entering the new data.
$ws.Activate()
$ws.Range("E36").Select()
